$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the existing hyperlink (on A7: trabajoAcademico) up front; it will be
#     (re)created, together with the other three RENATI type URIs, once the new
#     title row has been inserted and everything sits at its final address. ---
$ws.Range("A7").Hyperlinks.Delete()

# --- Insert a new first row for the "title" metadata entry ---
$ws.Rows("1:1").Insert()
$ws.Rows("1:1").RowHeight = 15.75

# Give A1 the same (bold) look as A2 ("namespace" label) without inventing a new style.
$ws.Range("A2").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A1").Value = "title"
$ws.Range("B1").Value = "Tipos de trabajo de investigación (RENATI)"

# --- Give A5:A8 (the RENATI type URI cells, now shifted down a row) one
#     consistent starting format before the hyperlinks are (re)applied, so
#     they end up sharing a single cell style. ---
$ws.Range("A5").Copy()
$ws.Range("A5:A8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Refresh the RENATI-type URI cells to the new "#fragment" form ---
$ws.Range("A5").Value = "http://purl.org/pe-repo/renati/tipo#tesis"
$ws.Range("A6").Value = "http://purl.org/pe-repo/renati/tipo#trabajoDeInvestigacion"
$ws.Range("A7").Value = "http://purl.org/pe-repo/renati/tipo#trabajoDeSuficienciaProfesional"
$ws.Range("A8").Value = "http://purl.org/pe-repo/renati/tipo#trabajoAcademico"

# --- (Re)create hyperlinks, anchored on the base namespace URI + in-page location ---
$ws.Hyperlinks.Add($ws.Range("A8"), "http://purl.org/pe-repo/renati/tipo", "trabajoAcademico")
$ws.Hyperlinks.Add($ws.Range("A5"), "http://purl.org/pe-repo/renati/tipo", "tesis")
$ws.Hyperlinks.Add($ws.Range("A6"), "http://purl.org/pe-repo/renati/tipo", "trabajoDeInvestigacion")
$ws.Hyperlinks.Add($ws.Range("A7"), "http://purl.org/pe-repo/renati/tipo", "trabajoDeSuficienciaProfesional")
